{"js": "// Replace the 25 \"dividend\u00f7divisor=\" expression strings in the practice-sheet\n// table with the newly generated set of problems. Each old expression is\n// unique in the document, so a plain text search-and-replace (matchCase, no\n// wildcards) targets exactly the one run that needs to change.\nconst replacements = [\n  [\"143\u00f79=\", \"669\u00f78=\"],\n  [\"823\u00f73=\", \"632\u00f79=\"],\n  [\"679\u00f74=\", \"363\u00f72=\"],\n  [\"683\u00f76=\", \"331\u00f78=\"],\n  [\"900\u00f72=\", \"819\u00f75=\"],\n  [\"509\u00f72=\", \"641\u00f78=\"],\n  [\"570\u00f73=\", \"860\u00f78=\"],\n  [\"657\u00f75=\", \"649\u00f78=\"],\n  [\"249\u00f73=\", \"129\u00f79=\"],\n  [\"655\u00f78=\", \"786\u00f72=\"],\n  [\"306\u00f72=\", \"517\u00f76=\"],\n  [\"391\u00f76=\", \"149\u00f72=\"],\n  [\"407\u00f78=\", \"471\u00f72=\"],\n  [\"764\u00f78=\", \"234\u00f74=\"],\n  [\"338\u00f76=\", \"470\u00f72=\"],\n  [\"960\u00f74=\", \"917\u00f75=\"],\n  [\"287\u00f76=\", \"460\u00f73=\"],\n  [\"225\u00f72=\", \"906\u00f78=\"],\n  [\"278\u00f72=\", \"301\u00f76=\"],\n  [\"620\u00f79=\", \"190\u00f79=\"],\n  [\"503\u00f76=\", \"810\u00f78=\"],\n  [\"541\u00f76=\", \"683\u00f72=\"],\n  [\"938\u00f72=\", \"511\u00f73=\"],\n  [\"272\u00f74=\", \"731\u00f78=\"],\n  [\"855\u00f75=\", \"589\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"dividend\u00f7divisor=\" expression strings in the practice-sheet\n# table with the newly generated set of problems. Each old expression is\n# unique in the document, so Find/Replace (MatchCase, whole-document range)\n# targets exactly the one run that needs to change per pair.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{old=\"143\u00f79=\"; new=\"669\u00f78=\"},\n    @{old=\"823\u00f73=\"; new=\"632\u00f79=\"},\n    @{old=\"679\u00f74=\"; new=\"363\u00f72=\"},\n    @{old=\"683\u00f76=\"; new=\"331\u00f78=\"},\n    @{old=\"900\u00f72=\"; new=\"819\u00f75=\"},\n    @{old=\"509\u00f72=\"; new=\"641\u00f78=\"},\n    @{old=\"570\u00f73=\"; new=\"860\u00f78=\"},\n    @{old=\"657\u00f75=\"; new=\"649\u00f78=\"},\n    @{old=\"249\u00f73=\"; new=\"129\u00f79=\"},\n    @{old=\"655\u00f78=\"; new=\"786\u00f72=\"},\n    @{old=\"306\u00f72=\"; new=\"517\u00f76=\"},\n    @{old=\"391\u00f76=\"; new=\"149\u00f72=\"},\n    @{old=\"407\u00f78=\"; new=\"471\u00f72=\"},\n    @{old=\"764\u00f78=\"; new=\"234\u00f74=\"},\n    @{old=\"338\u00f76=\"; new=\"470\u00f72=\"},\n    @{old=\"960\u00f74=\"; new=\"917\u00f75=\"},\n    @{old=\"287\u00f76=\"; new=\"460\u00f73=\"},\n    @{old=\"225\u00f72=\"; new=\"906\u00f78=\"},\n    @{old=\"278\u00f72=\"; new=\"301\u00f76=\"},\n    @{old=\"620\u00f79=\"; new=\"190\u00f79=\"},\n    @{old=\"503\u00f76=\"; new=\"810\u00f78=\"},\n    @{old=\"541\u00f76=\"; new=\"683\u00f72=\"},\n    @{old=\"938\u00f72=\"; new=\"511\u00f73=\"},\n    @{old=\"272\u00f74=\"; new=\"731\u00f78=\"},\n    @{old=\"855\u00f75=\"; new=\"589\u00f74=\"}\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.Text = $pair.old\n    $find.Replacement.Text = $pair.new\n    $find.Execute($pair.old, $false, $false, $false, $false, $false, $true, 1, $true, $pair.new, 2)\n}\n"}
